# Fix element ordering inside <w:rPr> of several Pandoc "highlighting"
# character styles in styles.xml: the <w:b/>/<w:i/> (bold/italic) toggle
# elements were serialized after <w:color>, which violates the CT_RPr
# sequence in wml.xsd (color must come before the bold/italic group).
# Re-touching each style's Font.Bold / Font.Italic (round-tripping the
# existing value) makes Word rewrite the run-properties in the
# schema-correct order without changing any actual formatting.

$d = $word.ActiveDocument

# Only the styles that actually carry a <w:b/> and/or <w:i/> toggle need
# to be touched, and only that same toggle should be re-applied -- round
# tripping the *other* (already-absent/default) one would turn an absent
# element into an explicit w:val="0", which is not part of this fix.
$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
$italicOnly = @("CommentTok", "DocumentationTok")
$boldAndItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")

foreach ($styleId in $boldOnly) {
    $style = $d.Styles.Item($styleId)
    $style.Font.Bold = $style.Font.Bold
}

foreach ($styleId in $italicOnly) {
    $style = $d.Styles.Item($styleId)
    $style.Font.Italic = $style.Font.Italic
}

foreach ($styleId in $boldAndItalic) {
    $style = $d.Styles.Item($styleId)
    $style.Font.Bold = $style.Font.Bold
    $style.Font.Italic = $style.Font.Italic
}
